# Generate Report for Handoff
# Updates the localization-status report: the "75cb5be8..." source file has
# been newly handed off for translation (status -> "Ready for handoff",
# priority -> "mt", new handoff timestamps/files), while the report rows are
# re-sorted by status so the "ed173f26..." (still "In Translation") file now
# sits in row 2 and "75cb5be8..." drops to row 3.

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db738e00bd107c5b220ab4222edb55d8d1d67fb3/e2e/75cb5be8-787a-4291-9980-2c2fb8563f68.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db738e00bd107c5b220ab4222edb55d8d1d67fb3/e2e/ed173f26-c992-400b-9b42-4e947fcbd5ba.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 2 now describes ed173f26 (unchanged status/date), row 3 now describes
# 75cb5be8 (status + date updated to reflect the new handoff).
$ws.Range("A2").Value = "ed173f26-c992-400b-9b42-4e947fcbd5ba.md"
$ws.Range("B2").Value = "e2e\ed173f26-c992-400b-9b42-4e947fcbd5ba.md"
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("G2").Value = "2016-08-31 12:15:07"

$ws.Range("A3").Value = "75cb5be8-787a-4291-9980-2c2fb8563f68.md"
$ws.Range("B3").Value = "e2e\75cb5be8-787a-4291-9980-2c2fb8563f68.md"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-31 12:15:56"

# Hyperlinks keep pointing at the same targets (rId2 -> 75cb5be8, rId3 ->
# ed173f26) but their display text now matches the swapped row order.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $urlA, "", "", "e2e\ed173f26-c992-400b-9b42-4e947fcbd5ba.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $urlB, "", "", "e2e\75cb5be8-787a-4291-9980-2c2fb8563f68.md")

# Columns E/F widened (status text got longer: "Ready for handoff").
$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ed173f26-c992-400b-9b42-4e947fcbd5ba.md"
$ws.Range("G2").Value = "ed173f26-c992-400b-9b42-4e947fcbd5ba.2c403bc7b7ae7045d8fa8f29d5def682bf243a9d.zh-cn.xlf"

$ws.Range("A3").Value = "75cb5be8-787a-4291-9980-2c2fb8563f68.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "75cb5be8-787a-4291-9980-2c2fb8563f68.a3ae6ce02551d12b210f87be10e626eef9b44e58.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-31 12:15:51"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlA, "", "", "ed173f26-c992-400b-9b42-4e947fcbd5ba.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlB, "", "", "75cb5be8-787a-4291-9980-2c2fb8563f68.md")

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ed173f26-c992-400b-9b42-4e947fcbd5ba.md"
$ws.Range("G2").Value = "ed173f26-c992-400b-9b42-4e947fcbd5ba.2c403bc7b7ae7045d8fa8f29d5def682bf243a9d.de-de.xlf"

$ws.Range("A3").Value = "75cb5be8-787a-4291-9980-2c2fb8563f68.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "75cb5be8-787a-4291-9980-2c2fb8563f68.a3ae6ce02551d12b210f87be10e626eef9b44e58.de-de.xlf"
$ws.Range("H3").Value = "2016-08-31 12:15:56"

$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlA, "", "", "ed173f26-c992-400b-9b42-4e947fcbd5ba.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlB, "", "", "75cb5be8-787a-4291-9980-2c2fb8563f68.md")

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332
